$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Main Data")
$ws2 = $wb.Worksheets.Item("Exporatable")

# --- Fill in the "Non PCA" (Normal NN (large)) performance row with data ---
# Left block (In Sample), columns B:G
$ws1.Range("B9").Value = 0.97791681687037102
$ws1.Range("C9").Value = 0.98050443331400505
$ws1.Range("D9").Formula = "=(0.981132060289382 + 0.978257079919179)/2"
$ws1.Range("E9").Value = 0.98355424404144198
$ws1.Range("F9").Formula = "=(0.97914108633995 + 0.990791896979014)/2"
$ws1.Range("G9").Value = 0.99015988906224495

# Right block (Out of Sample), columns K:P
$ws1.Range("K9").Value = 0.90392156442006399
$ws1.Range("L9").Value = 0.90686275561650598
$ws1.Range("M9").Formula = "=(0.911764681339263 + 0.901960770289103)/2"
$ws1.Range("N9").Value = 0.89411763350168805
$ws1.Range("O9").Formula = "=(0.907352954149246 + 0.877450982729593) /2"
$ws1.Range("P9").Value = 0.87941175699233998

# Credit note next to the row, styled with a left/center aligned black font
$ws1.Range("R9").Value = "Jesse"

# M9 picks up a distinct font + alignment (no borders) rather than the
# table's usual bordered/centered numeric style
$ws1.Range("M9").Font.Color = 0
$ws1.Range("M9").HorizontalAlignment = -4131
$ws1.Range("M9").VerticalAlignment = -4108
$ws1.Range("M9").Borders.LineStyle = -4142

# --- View/selection state ---
$ws2.Range("B2").Select()
$ws1.Activate()
$ws1.Range("Q23").Select()
